# model_param_lookup.xlsx update
# - adds a new "model_description" column (H) describing each experiment
# - adds a legend row (row 2: ours.png / base(deep_voxel_flow))
# - re-orders a few experiment rows and corrects several coefficient values
# - appends two new experiment rows (15, 16) logged on 2019-04-23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "model_id"
$ws.Range("B1").Value = "strategy"
$ws.Range("C1").Value = "coef_loss_e"
$ws.Range("D1").Value = "coef_loss_t"
$ws.Range("E1").Value = "coef_loss_m"
$ws.Range("F1").Value = "s1_epochs"
$ws.Range("G1").Value = "max_epochs"
$ws.Range("H1").Value = "model_description"

# Row 2
$ws.Range("A2").Value = "ours.png"
$ws.Range("B2").Value = "base(deep_voxel_flow)"
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").Value = "[Baseline]  Deep Voxel Flow"

# Row 3
$ws.Range("A3").Value = "2019-04-14T232535"
$ws.Range("B3").Value = "original_cycle_gen"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.1
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = "[Proposed] CyclicGen"

# Row 4
$ws.Range("A4").Value = "2019-04-15T111748"
$ws.Range("B4").Value = "original_cycle_gen"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = "[Proposed] CyclicGen"

# Row 5
$ws.Range("A5").Value = "2019-04-17T191452"
$ws.Range("B5").Value = "original_cycle_gen"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "[Baseline]  1 Batch Only"

# Row 6
$ws.Range("A6").Value = "2019-04-18T164451"
$ws.Range("B6").Value = "accel_adjust"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = "[Proposed] CyclicGen: Adapt to Acceleration "

# Row 7
$ws.Range("A7").Value = "2019-04-18T233637"
$ws.Range("B7").Value = "accel_adjust"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.1
$ws.Range("E7").Value = 0.1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = "[Proposed] CyclicGen: Adapt to Acceleration "

# Row 8
$ws.Range("A8").Value = "2019-04-19T174028"
$ws.Range("B8").Value = "original_cycle_gen"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.1
$ws.Range("F8").Value = 5
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = "[Baseline]  Stage 1 Only"

# Row 9
$ws.Range("A9").Value = "2019-04-20T012234"
$ws.Range("B9").Value = "original_cycle_gen"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.1
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 10
$ws.Range("H9").Value = "[Baseline]  CyclicGen"

# Row 10
$ws.Range("A10").Value = "2019-04-20T151155"
$ws.Range("B10").Value = "accel_adjust"
$ws.Range("C10").Value = 0.5
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.1
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 10
$ws.Range("H10").Value = "[Proposed] CyclicGen: Adapt to Acceleration & Extra Cycle Consistency Loss (coef: 0.5)"

# Row 11
$ws.Range("A11").Value = "2019-04-20T215827"
$ws.Range("B11").Value = "accel_adjust"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.1
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = "[Proposed] CyclicGen: Adapt to Acceleration "

# Row 12
$ws.Range("A12").Value = "2019-04-21T092001"
$ws.Range("B12").Value = "original_cycle_gen"
$ws.Range("C12").Value = 0.5
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.1
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = "[Proposed] CyclicGen: Extra Cycle Consistency Loss (coef: 0.5)"

# Row 13
$ws.Range("A13").Value = "2019-04-21T224636"
$ws.Range("B13").Value = "accel_adjust"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.1
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = "[Proposed] CyclicGen: Adapt to Acceleration & Extra Cycle Consistency Loss (coef: 1.0)"

# Row 14
$ws.Range("A14").Value = "2019-04-22T082112"
$ws.Range("B14").Value = "original_cycle_gen"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.1
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 10
$ws.Range("H14").Value = "[Proposed] CyclicGen: Extra Cycle Consistency Loss (coef: 1.0)"

# Row 15
$ws.Range("A15").Value = "2019-04-23T083858"
$ws.Range("B15").Value = "original_cycle_gen"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = "[Proposed] CyclicGen: Extra Cycle Consistency Loss (coef: 1.0)"

# Row 16
$ws.Range("A16").Value = "2019-04-23T141504"
$ws.Range("B16").Value = "accel_adjust"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = "[Proposed] CyclicGen: Adapt to Acceleration "

# Column widths to fit the new/changed columns
$ws.Columns.Item(2).ColumnWidth = 20.46484375
$ws.Columns.Item(3).ColumnWidth = 9.796875
$ws.Columns.Item(5).ColumnWidth = 11.1328125
$ws.Columns.Item(7).ColumnWidth = 12
$ws.Columns.Item(8).ColumnWidth = 55.265625

# Restore the view: scrolled one column right, H6 selected
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("H6").Select()

